# This script applies the recalculated model results from the commit
# "additional F_min expressions tested": component names are renamed
# (PLP -> Ald, T3H -> Hydr, Comp -> Hydrz) and the dependent numeric
# results throughout the workbook are updated to the new fit values.

$wb = $excel.ActiveWorkbook

# ---- 1) Rename component name labels (text, kept as text) ----

$ws = $wb.Worksheets.Item("input_stoich_coefficients")
$ws.Range("A1").Value = "Ald"
$ws.Range("B1").Value = "Hydr"

$ws = $wb.Worksheets.Item("input_concentrations")
$ws.Range("A2").Value = "Ald"
$ws.Range("B2").Value = "Hydr"

$ws = $wb.Worksheets.Item("equilibrium_concentrations")
$ws.Range("A1").Value = "Ald"
$ws.Range("B1").Value = "Hydr"
$ws.Range("C1").Value = "Hydrz"

$ws = $wb.Worksheets.Item("component_names")
$ws.Range("A1").Value = "Hydr"

$ws = $wb.Worksheets.Item("constants_evaluated")
$ws.Range("A1").Value = "Ald"
$ws.Range("A2").Value = "Hydr"
$ws.Range("A3").Value = "Hydrz"

$ws = $wb.Worksheets.Item("correlation_matrix")
$ws.Range("B1").Value = "Hydrz"
$ws.Range("A2").Value = "Hydrz"

$ws = $wb.Worksheets.Item("enthalpies_calc")
$ws.Range("B1").Value = "Ald"
$ws.Range("C1").Value = "Hydr"
$ws.Range("D1").Value = "Hydrz"

# ---- 2) Text cells holding numeric-looking evaluated-constant readouts ----
# (leading apostrophe forces these to stay text, matching the original type)
$ws = $wb.Worksheets.Item("constants_evaluated")
$ws.Range("B3").Value = "'5.1243896484375"
$ws.Range("C3").Value = "'0.08683091631657618"

# ---- 3) Updated numeric results (recalculated model outputs) ----

$ws = $wb.Worksheets.Item("enthalpies_calc")
$ws.Range("B2").Value = [double]"-0"
$ws.Range("C2").Value = [double]"-0"
$ws.Range("D2").Value = [double]"-47.18400184222095"
$ws.Range("D3").Value = [double]"0.5878688796803043"

$ws = $wb.Worksheets.Item("equilibrium_concentrations")
$ws.Range("B2").Value = [double]"9.333989435922405e-21"
$ws.Range("C2").Value = [double]"9.906660105643585e-19"
$ws.Range("A3").Value = [double]"0.0006976287261882581"
$ws.Range("B3").Value = [double]"1.055626040551382e-06"
$ws.Range("C3").Value = [double]"9.806727295240044e-05"
$ws.Range("A4").Value = [double]"0.0005989087469879397"
$ws.Range("B4").Value = [double]"2.450869330314637e-06"
$ws.Range("C4").Value = [double]"0.0001954656560849638"
$ws.Range("A5").Value = [double]"0.0005010511050819636"
$ws.Range("B5").Value = [double]"4.37643244036896e-06"
$ws.Range("C5").Value = [double]"0.0002920060848030009"
$ws.Range("A6").Value = [double]"0.0004044140977156763"
$ws.Range("B6").Value = [double]"7.19226151289438e-06"
$ws.Range("C6").Value = [double]"0.0003873302400944918"
$ws.Range("A7").Value = [double]"0.0003097512590243261"
$ws.Range("B7").Value = [double]"1.165352867947256e-05"
$ws.Range("C7").Value = [double]"0.0004806845662014034"
$ws.Range("A8").Value = [double]"0.0002188669363957195"
$ws.Range("B8").Value = [double]"1.956620851367137e-05"
$ws.Range("C8").Value = [double]"0.0005702646942558629"
$ws.Range("A9").Value = [double]"0.0001366229432178645"
$ws.Range("B9").Value = [double]"3.579373086772801e-05"
$ws.Range("C9").Value = [double]"0.0006512087895313295"
$ws.Range("A10").Value = [double]"7.447759817241203e-05"
$ws.Range("B10").Value = [double]"7.179602025413232e-05"
$ws.Range("C10").Value = [double]"0.0007120585121480837"
$ws.Range("A11").Value = [double]"4.103791556752228e-05"
$ws.Range("B11").Value = [double]"0.0001361816862657428"
$ws.Range("C11").Value = [double]"0.0007442068267402886"
$ws.Range("A12").Value = [double]"2.602781679037327e-05"
$ws.Range("B12").Value = [double]"0.0002186762351147522"
$ws.Range("C12").Value = [double]"0.0007579297910064992"
$ws.Range("A13").Value = [double]"1.860280597664169e-05"
$ws.Range("B13").Value = [double]"0.0003084367453874388"
$ws.Range("C13").Value = [double]"0.0007640718800534368"
$ws.Range("A14").Value = [double]"1.436188864504015e-05"
$ws.Range("B14").Value = [double]"0.0004010637867631738"
$ws.Range("C14").Value = [double]"0.0007670340677972558"
$ws.Range("A15").Value = [double]"1.166010990015545e-05"
$ws.Range("B15").Value = [double]"0.000494913958300624"
$ws.Range("C15").Value = [double]"0.0007684612888361514"
$ws.Range("A16").Value = [double]"9.800205155534634e-06"
$ws.Range("B16").Value = [double]"0.0005892915392431673"
$ws.Range("C16").Value = [double]"0.0007690507878610797"
$ws.Range("A17").Value = [double]"8.445915395269162e-06"
$ws.Range("B17").Value = [double]"0.0006838618043618297"
$ws.Range("C17").Value = [double]"0.0007691388046014878"

$ws = $wb.Worksheets.Item("heats_calc_abs_errors")
$ws.Range("B3").Value = [double]"0.002380868451341722"
$ws.Range("C3").Value = [double]"-0.002237855946052733"
$ws.Range("D3").Value = [double]"0.001131751355276961"
$ws.Range("E3").Value = [double]"0.002657222003472431"
$ws.Range("F3").Value = [double]"-0.0006787621323298826"
$ws.Range("G3").Value = [double]"-0.0001663928370761597"
$ws.Range("H3").Value = [double]"-0.003683889206730295"
$ws.Range("I3").Value = [double]"-0.0003882809624511721"
$ws.Range("J3").Value = [double]"0.002885421085913706"
$ws.Range("K3").Value = [double]"-0.001784150505066391"
$ws.Range("L3").Value = [double]"-0.00343140357250874"
$ws.Range("M3").Value = [double]"-0.0007850088946899834"
$ws.Range("N3").Value = [double]"-0.001947310028463187"
$ws.Range("O3").Value = [double]"-0.002562551125734338"
$ws.Range("P3").Value = [double]"-0.00246694741233297"

$ws = $wb.Worksheets.Item("heats_calc_rel_errors")
$ws.Range("B3").Value = [double]"0.03311151865466019"
$ws.Range("C3").Value = [double]"0.03337992473480132"
$ws.Range("D3").Value = [double]"0.01616074499805673"
$ws.Range("E3").Value = [double]"0.03746928154181674"
$ws.Range("F3").Value = [double]"0.01022227970352956"
$ws.Range("G3").Value = [double]"0.002582345720644832"
$ws.Range("H3").Value = [double]"0.06704693846773802"
$ws.Range("I3").Value = [double]"0.008820582429855261"
$ws.Range("J3").Value = [double]"0.1075892118464986"
$ws.Range("K3").Value = [double]"0.1989640362769153"
$ws.Range("L3").Value = [double]"1.81692542719134"
$ws.Range("M3").Value = [double]"0.3482936222569896"
$ws.Range("N3").Value = [double]"-171.5238288085256"
$ws.Range("O3").Value = [double]"-2.083446041030945"
$ws.Range("P3").Value = [double]"-1.648192803472683"

$ws = $wb.Worksheets.Item("correlation_matrix")
$ws.Range("B2").Value = [double]"1"
